$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The row representing requirement T4.2 (old placeholder "X"/"X") has been
# implemented in software ("latch implemented in software"), so the whole
# row is removed. This shifts all subsequent rows up by one.
$deletedRowHeight = $ws.Rows(12).RowHeight
$ws.Rows(12).Delete()

# The pictures below the deleted row are anchored to cells but this runtime
# does not automatically re-anchor them when a row is removed, so nudge them
# up by the height of the removed row to keep them aligned with their rows.
foreach ($sh in $ws.Shapes) {
    $sh.Top = $sh.Top - $deletedRowHeight
}

# The requirement that used to be T4.3 is now renumbered to T4.2, since the
# old T4.2 placeholder row no longer exists.
$ws.Range("A12").Value = "T4.2"

# Reflect the cursor/selection position left behind in the saved file.
$ws.Range("C11").Select() | Out-Null

# The "_ftn1" footnote anchor pointed at A14; after the row deletion above,
# that content now lives at A13, so repoint the defined name accordingly.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Blad1!_ftn1") {
        $n.RefersTo = "=Blad1!`$A`$13"
    }
}
